$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: merge the two runs "A" + " probléma bemutatása…"
#    into a single run with identical formatting.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "A probléma bemutatása és erre a javasolt megoldás kifejtése szövegesen",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "A probléma bemutatása és erre a javasolt megoldás kifejtése szövegesen",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the stale _GoBack bookmark from the end of the "Mindezt…"
#    paragraph - it gets relocated to the start of the rewritten
#    "Ezután…" paragraph below.
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 3) Rewrite the "Ezután … eljut az érme a tárolóba." paragraph with
#    the new, much longer description (and move the _GoBack bookmark
#    to the top of it).
# ------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$r7 = $p7.Range

$xml7 = "<w:p $wns>" +
  "<!-- placeholder -->" +
  "<w:pPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
  "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr>" +
  "<w:t xml:space=`"preserve`">Ezután a sínre merőleges lemez ami eddig biztosította, hogy az érme ne tudjon tovább gurulni, hanem fix állapotba legyen, most egy léptető motor segítségével elfordul a motor tengelyével párhuzamosan amíg az érem legurul /leesik a sínvégéről, majd vissza kerül az eredeti helyére. Mikor az érme </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>leesik</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr>" +
  "<w:t xml:space=`"preserve`"> a sínről egy  lapos alacsony perű fém lapra esik melynek szélessége nagyobb mint a kétszázas átmérője. Ennek a lapnak az aljára egy cső van rögzítve melynek végén egy </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>servo</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr>" +
  "<w:t xml:space=`"preserve`"> motor van, ami bizonyos helyzetekbe állítja a fémlapot, hogy a vége a megfelelő tárolóedény fölött legyen, így az érmek – mely ekkor már nem </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>élén</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr>" +
  "<w:t xml:space=`"preserve`"> hanem lapján fekszik – az enyhén döntött lemezen le tud csúszni a helyére. </w:t></w:r>" +
  "</w:p>"

$r7.InsertXML($xml7)

# ------------------------------------------------------------------
# 4) Rewrite the "Ahonnan … segítségével ki kerül az érme." paragraph,
#    adding the <w:lastRenderedPageBreak/> marker now that the
#    preceding paragraph pushed it onto a new page.
# ------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$r8 = $p8.Range

$xml8 = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr>" +
  "<w:lastRenderedPageBreak/><w:t>Ahonnan …</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr>" +
  "<w:t xml:space=`"preserve`"> segítségével ki kerül az érme.</w:t></w:r>" +
  "</w:p>"

$r8.InsertXML($xml8)
